{"js": "// Office.js (Word JavaScript API) edit script.\n// Rewrites the \"history\" themed document into the \"mathematics\" themed\n// document described by the commit diff: title, author, email, the two\n// long body paragraphs, and a trailing empty paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Paragraph 0: Title\n// ---------------------------------------------------------------------\nparagraphs.items[0].insertText(\n  \"The Mathematical Keys: Unlocking the Enigma of the Universe\",\n  Word.InsertLocation.replace\n);\n\n// ---------------------------------------------------------------------\n// Paragraph 1: Author name\n// ---------------------------------------------------------------------\nparagraphs.items[1].insertText(\n  \"Richard Albert\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Paragraph 2: Email address, built from several runs: \"jane\" + \".\" +\n// \"longfellow@historian\" + \".\" + \"edu\". Only the two non-literal \".\"\n// runs stay the same, so target the two text runs individually via\n// search (keeps the existing \".\" runs and their formatting intact).\n// ---------------------------------------------------------------------\nconst userRange = body.search(\"jane\", { matchCase: true, matchWholeWord: true });\nuserRange.load(\"text\");\nconst domainRange = body.search(\"longfellow@historian\", { matchCase: true });\ndomainRange.load(\"text\");\nawait context.sync();\n\nuserRange.items[0].insertText(\"Richard\", Word.InsertLocation.replace);\ndomainRange.items[0].insertText(\"Albert@CrestPoint\", Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Paragraph 4: First long body paragraph (three sub-blocks separated by\n// double manual line breaks). Replace the whole paragraph text in one\n// shot -- the run formatting (Calibri, black, sz 24) carries over from\n// the paragraph's existing run.\n// ---------------------------------------------------------------------\nconst bodyPara1 =\n  \"Within the vast cosmic tapestry, mathematics serves as a powerful beacon, illuminating the enigmatic mysteries that surround us.\" +\n  \" It is a universal language that transcends cultures and time, providing a framework for comprehending the intricate workings of the cosmos.\" +\n  \" Through the study of mathematics, we unravel the symphony of patterns hidden in nature, unveiling the fundamental principles underpinning the universe.\" +\n  \" Like a quantum dance of numbers, mathematical concepts orchestrate the movement of the stars, the growth of organisms, and the flow of time itself.\" +\n  \" As we delve into the depths of mathematics, we uncover the secrets of the universe, revealing a mesmerizing enigma that captivates and inspires.\" +\n  \" This journey into the realm of mathematics unveils the profound interconnectedness between humanity and the cosmos, forever expanding our understanding of the universe and our place within it.\" +\n  \"\\v\\v\" +\n  \"In the realm of mathematics, we unlock the secrets to comprehend the world's fundamental building blocks and their intricate relationships.\" +\n  \" Abstract concepts like numbers, shapes, and equations unveil the underlying patterns that govern the universe.\" +\n  \" Whether exploring the infinitesimally small subatomic particles or the vastness of galaxies, mathematics provides a lens through which we discover the profound interconnectedness of all things.\" +\n  \" By unraveling the enigma of numbers, we embark on a quest to unveil the fundamental principles that orchestrate the symphony of existence, revealing the exquisite beauty and elegance of the universe.\" +\n  \"\\v\\v\" +\n  \"Mathematics, with its intrinsic beauty and timeless truths, has captivated the hearts and minds of scholars throughout history.\" +\n  \" From the ancient Egyptians and Babylonians to the brilliance of Isaac Newton and Albert Einstein, mathematics has served as a bridge connecting diverse cultures and civilizations.\" +\n  \" It has propelled technological advancements, facilitated scientific discoveries, and empowered us to navigate the complexities of our world.\" +\n  \" As we continue to push the boundaries of mathematical knowledge, we unlock new vistas of understanding, opening doors to realms beyond our current comprehension.\" +\n  \" This enduring pursuit of mathematical enlightenment enriches our lives, broadens our perspectives, and inspires us to embrace the wonders of the universe.\";\n\nparagraphs.items[4].insertText(bodyPara1, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Paragraph 6: Summary body paragraph.\n// ---------------------------------------------------------------------\nconst bodyPara2 =\n  \"Mathematics, with its enigmatic beauty and universal language, empowers us to unlock the mysteries of the universe.\" +\n  \" Through the study of numbers, shapes, and equations, we uncover the fundamental principles orchestrating the symphony of existence.\" +\n  \" Its enduring pursuit enriches our lives, broadens our perspectives, and propels us toward a deeper understanding of our world.\" +\n  \" Mathematics, a beacon of enlightenment, continues to reveal the interconnectedness of all things, igniting our imagination and inspiring our quest for knowledge.\";\n\nparagraphs.items[6].insertText(bodyPara2, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Append a new, empty paragraph at the very end of the document body.\n// ---------------------------------------------------------------------\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Rewrites the \"history\" themed document into the \"mathematics\" themed\n# document described by the commit diff: title, author, email, the two\n# long body paragraphs, and a trailing empty paragraph.\n\n$d = $word.ActiveDocument\n\n# -----------------------------------------------------------------\n# Paragraph 1: Title\n# -----------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Melody of Time: Echoes of History Resonating Today\"\n$find.Replacement.Text = \"The Mathematical Keys: Unlocking the Enigma of the Universe\"\n$find.MatchCase = $true\n$find.Forward = $true\n$find.Wrap = 0\n$find.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2) | Out-Null\n\n# -----------------------------------------------------------------\n# Paragraph 2: Author name\n# -----------------------------------------------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"Jane Longfellow\"\n$find2.Replacement.Text = \"Richard Albert\"\n$find2.MatchCase = $true\n$find2.Forward = $true\n$find2.Wrap = 0\n$find2.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2) | Out-Null\n\n# -----------------------------------------------------------------\n# Paragraph 3: Email address -- built from several runs:\n#   \"jane\" + \".\" + \"longfellow@historian\" + \".\" + \"edu\"\n# Only the two text runs (\"jane\" and \"longfellow@historian\") change;\n# the \".\" runs in between stay as-is.\n# -----------------------------------------------------------------\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Text = \"jane\"\n$find3.Replacement.Text = \"Richard\"\n$find3.MatchCase = $true\n$find3.MatchWholeWord = $true\n$find3.Forward = $true\n$find3.Wrap = 0\n$find3.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2) | Out-Null\n\n$find4 = $d.Content.Find\n$find4.ClearFormatting()\n$find4.Text = \"longfellow@historian\"\n$find4.Replacement.Text = \"Albert@CrestPoint\"\n$find4.MatchCase = $true\n$find4.MatchWholeWord = $false\n$find4.Forward = $true\n$find4.Wrap = 0\n$find4.Execute([ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,[ref]$null,2) | Out-Null\n\n# -----------------------------------------------------------------\n# Paragraph 5 (index 5 in the Paragraphs collection, 1-based):\n# first long body paragraph -- three blocks separated by double manual\n# line breaks (chr 11 = vertical-tab = Word's manual line break).\n# Replace the whole paragraph range in one shot so the run's existing\n# formatting (Calibri, black, sz 24) carries over intact.\n# -----------------------------------------------------------------\n$br = [char]11\n$bodyPara1 = \"Within the vast cosmic tapestry, mathematics serves as a powerful beacon, illuminating the enigmatic mysteries that surround us.\" +\n  \" It is a universal language that transcends cultures and time, providing a framework for comprehending the intricate workings of the cosmos.\" +\n  \" Through the study of mathematics, we unravel the symphony of patterns hidden in nature, unveiling the fundamental principles underpinning the universe.\" +\n  \" Like a quantum dance of numbers, mathematical concepts orchestrate the movement of the stars, the growth of organisms, and the flow of time itself.\" +\n  \" As we delve into the depths of mathematics, we uncover the secrets of the universe, revealing a mesmerizing enigma that captivates and inspires.\" +\n  \" This journey into the realm of mathematics unveils the profound interconnectedness between humanity and the cosmos, forever expanding our understanding of the universe and our place within it.\" +\n  \"$br$br\" +\n  \"In the realm of mathematics, we unlock the secrets to comprehend the world's fundamental building blocks and their intricate relationships.\" +\n  \" Abstract concepts like numbers, shapes, and equations unveil the underlying patterns that govern the universe.\" +\n  \" Whether exploring the infinitesimally small subatomic particles or the vastness of galaxies, mathematics provides a lens through which we discover the profound interconnectedness of all things.\" +\n  \" By unraveling the enigma of numbers, we embark on a quest to unveil the fundamental principles that orchestrate the symphony of existence, revealing the exquisite beauty and elegance of the universe.\" +\n  \"$br$br\" +\n  \"Mathematics, with its intrinsic beauty and timeless truths, has captivated the hearts and minds of scholars throughout history.\" +\n  \" From the ancient Egyptians and Babylonians to the brilliance of Isaac Newton and Albert Einstein, mathematics has served as a bridge connecting diverse cultures and civilizations.\" +\n  \" It has propelled technological advancements, facilitated scientific discoveries, and empowered us to navigate the complexities of our world.\" +\n  \" As we continue to push the boundaries of mathematical knowledge, we unlock new vistas of understanding, opening doors to realms beyond our current comprehension.\" +\n  \" This enduring pursuit of mathematical enlightenment enriches our lives, broadens our perspectives, and inspires us to embrace the wonders of the universe.\"\n\n$p1 = $d.Paragraphs.Item(5)\n$r1 = $d.Range($p1.Range.Start, $p1.Range.End)\n$r1.Text = $bodyPara1\n\n# -----------------------------------------------------------------\n# Paragraph 7 (1-based): Summary body paragraph.\n# -----------------------------------------------------------------\n$bodyPara2 = \"Mathematics, with its enigmatic beauty and universal language, empowers us to unlock the mysteries of the universe.\" +\n  \" Through the study of numbers, shapes, and equations, we uncover the fundamental principles orchestrating the symphony of existence.\" +\n  \" Its enduring pursuit enriches our lives, broadens our perspectives, and propels us toward a deeper understanding of our world.\" +\n  \" Mathematics, a beacon of enlightenment, continues to reveal the interconnectedness of all things, igniting our imagination and inspiring our quest for knowledge.\"\n\n$p2 = $d.Paragraphs.Item(7)\n$r2 = $d.Range($p2.Range.Start, $p2.Range.End)\n$r2.Text = $bodyPara2\n\n# -----------------------------------------------------------------\n# Append a new, empty paragraph at the very end of the document body.\n# -----------------------------------------------------------------\n$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null\n"}
